$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 13-15 belong to a new monster group "demo怪物组3" (id 3),
# following the same layout pattern as the existing groups 1 (rows 7-9)
# and 2 (rows 10-12).

# Seed formatting for the new rows by copying from the last existing
# data row (row 9), which already carries the correct cell styles
# for columns C-I (style ids 1/1/2/1/1/1/1).
$ws.Range("C9:I9").Copy() | Out-Null
$ws.Range("C13:I13").PasteSpecial(-4122) | Out-Null
$ws.Range("C14:I14").PasteSpecial(-4122) | Out-Null
$ws.Range("C15:I15").PasteSpecial(-4122) | Out-Null

# The first row of each group highlights column I (style id 5); mirror
# that for the first row of the new group too, copying from I10 (the
# first row of group 2) which already uses that highlighted style.
$ws.Range("I10").Copy() | Out-Null
$ws.Range("I13").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Row 13
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = "demo怪物组3"
$ws.Range("F13").Value = 31
$ws.Range("G13").Value = 32
$ws.Range("H13").Value = 30
$ws.Range("I13").Value = 1234

# Row 14
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = "demo怪物组3"
$ws.Range("F14").Value = 33
$ws.Range("G14").Value = 34
$ws.Range("H14").Value = 31
$ws.Range("I14").Value = 442

# Row 15
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = "demo怪物组3"
$ws.Range("F15").Value = 35
$ws.Range("G15").Value = 37
$ws.Range("H15").Value = 33
$ws.Range("I15").Value = 55194

# Match the saved selection state shown in the authored workbook.
$ws.Range("I13:I15").Select() | Out-Null
